$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 31

# Text columns: leading apostrophe keeps them as literal text instead of
# Excel auto-converting "2024-01-07" to a date serial or "01" to the number 1.
$ws.Cells.Item($row, 1).Value = "'2024-01-07"
$ws.Cells.Item($row, 2).Value = "21:03:42"
$ws.Cells.Item($row, 3).Value = "Sunday"
$ws.Cells.Item($row, 4).Value = "'01"

# Numeric columns
$ws.Cells.Item($row, 5).Value = 140554
$ws.Cells.Item($row, 6).Value = 143034
$ws.Cells.Item($row, 7).Value = 172327
$ws.Cells.Item($row, 8).Value = 147203
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 118405
$ws.Cells.Item($row, 11).Value = 224803
$ws.Cells.Item($row, 12).Value = 249629
$ws.Cells.Item($row, 13).Value = 185287
$ws.Cells.Item($row, 14).Value = 110454
$ws.Cells.Item($row, 15).Value = 40678
$ws.Cells.Item($row, 16).Value = 30808
$ws.Cells.Item($row, 17).Value = 72553
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42546
$ws.Cells.Item($row, 20).Value = -1
